$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.71
$ws.Range("G2").Value = 1.72
$ws.Range("L2").Value = 1.35
$ws.Range("Y2").Value = 21
$ws.Range("AB2").Value = 9.4
$ws.Range("AN2").Value = 9.800000000000001

# Row 3
$ws.Range("AH3").Value = 48
$ws.Range("AK3").Value = 22
$ws.Range("AL3").Value = 1000
$ws.Range("AN3").Value = 10.5

# Row 4
$ws.Range("F4").Value = 8.800000000000001
$ws.Range("J4").Value = 5.1
$ws.Range("L4").Value = 1.4
$ws.Range("N4").Value = 4.1
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 2.08
$ws.Range("R4").Value = 1.4
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 1.76
$ws.Range("V4").Value = 3.25
$ws.Range("Y4").Value = 8
$ws.Range("Z4").Value = 7.8
$ws.Range("AO4").Value = 7
